# Added Hawkeye for both configs using bzip2 benchmark
$wb = $excel.ActiveWorkbook

# ---- Config1 sheet (Hawkeye row 9, OPTGen row 10 for bzip2) ----
$ws1 = $wb.Worksheets.Item("Config1")

# Row 9 - Hawkeye
$ws1.Range("C9").Value = 50000001
$ws1.Range("D9").Value = 68593474
$ws1.Range("E9").Value = 620418
$ws1.Range("F9").Value = 560367
$ws1.Range("G9").Value = 60051
$ws1.Range("H9").Formula = "=(C9/D9)"
$ws1.Range("I9").Formula = "=G9/(C9/1000)"

# Row 10 - OPTGen
$ws1.Range("C10").Value = 50000001
$ws1.Range("D10").Value = 68593474
$ws1.Range("E10").Value = 12097
$ws1.Range("F10").Value = 10446
$ws1.Range("G10").Formula = "=E10-F10"
$ws1.Range("H10").Formula = "=(C10/D10)"
$ws1.Range("I10").Formula = "=G10/(C10/1000)"
$ws1.Range("J10").Formula = "=F10/E10"

$ws1.Range("C11").Select()

# ---- Config2 sheet (Hawkeye row 9, OPTGen row 10 for bzip2) ----
$ws2 = $wb.Worksheets.Item("Config2")

# Row 9 - Hawkeye
$ws2.Range("C9").Value = 50000001
$ws2.Range("D9").Value = 66883100
$ws2.Range("E9").Value = 853878
$ws2.Range("F9").Value = 787833
$ws2.Range("G9").Value = 66045
$ws2.Range("H9").Formula = "=(C9/D9)"
$ws2.Range("I9").Formula = "=G9/(C9/1000)"

# Row 10 - OPTGen
$ws2.Range("C10").Value = 50000001
$ws2.Range("D10").Value = 66883100
$ws2.Range("E10").Value = 10725
$ws2.Range("F10").Value = 12081
$ws2.Range("G10").Formula = "=E10-F10"
$ws2.Range("H10").Formula = "=(C10/D10)"
$ws2.Range("I10").Formula = "=G10/(C10/1000)"
$ws2.Range("J10").Formula = "=F10/E10"

$ws2.Range("C15").Select()

# Restore Config1 as the active/selected sheet (matches original tab state)
$ws1.Activate()
$ws1.Range("C11").Select()
